$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.851.88"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "2.278.11"
$ws.Range("E3").Value = "  +4.77%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'250.03"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "'0.636"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("D7").Value = "'72.14"
$ws.Range("E7").Value = "  +10.15%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.661"
$ws.Range("E9").Value = "  +17.62%  "
$ws.Range("D10").Value = "'38.84"
$ws.Range("E10").Value = "  +9.36%  "
$ws.Range("D11").Value = "'59.82"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'0.0970"
$ws.Range("E12").Value = "  +4.99%  "
$ws.Range("E13").Value = "  +8.54%  "
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "2.615.38"
$ws.Range("E15").Value = "  +4.54%  "
$ws.Range("D16").Value = "'14.87"
$ws.Range("E16").Value = "  +4.22%  "
$ws.Range("D17").Value = "'0.888"
$ws.Range("E17").Value = "  +5.03%  "
$ws.Range("D18").Value = "2.269.25"
$ws.Range("E18").Value = "  +4.23%  "
$ws.Range("D19").Value = "42.787.01"
$ws.Range("E19").Value = "  +4.01%  "
$ws.Range("D20").Value = "'0.0000101"
$ws.Range("E20").Value = "  +7.59%  "
$ws.Range("D21").Value = "'6.32"
$ws.Range("E21").Value = "  +4.05%  "
$ws.Range("D22").Value = "'73.19"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("D23").Value = "'235.81"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("D24").Value = "'2.15"
$ws.Range("E24").Value = "  +5.29%  "
$ws.Range("D25").Value = "'4.04"
$ws.Range("E25").Value = "  +4.64%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'11.39"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "'2.44"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("D30").Value = "'2.14"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").Value = "'167.61"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "'21.03"
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("D33").Value = "'6.43"
$ws.Range("E33").Value = "  +13.44%  "
$ws.Range("E34").Value = "  +4.85%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'31.78"
$ws.Range("E35").Value = "  +32.11%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0806"
$ws.Range("E36").Value = "  +9.16%  "
$ws.Range("E37").Value = "  +4.20%  "
$ws.Range("D38").Value = "'4.48"
$ws.Range("E38").Value = "  +13.69%  "
$ws.Range("D39").Value = "'4.78"
$ws.Range("E39").Value = "  +5.87%  "
$ws.Range("D40").Value = "'0.0314"
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("D41").Value = "'2.35"
$ws.Range("E41").Value = "  +7.51%  "
$ws.Range("D42").Value = "'12.80"
$ws.Range("E42").Value = "  +15.61%  "
$ws.Range("E43").Value = "  +6.91%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'5.08"
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.208"
$ws.Range("E45").Value = "  +9.76%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'9.29"
$ws.Range("E46").Value = "  +9.72%  "
$ws.Range("D47").Value = "'62.20"
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("E48").Value = "  +3.38%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'1.17"
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("E51").Value = "  +4.70%  "
